# Add the API tests for Sign Up & Assertions Enhancements
#
# - Rename "testsheet1" -> "API" and "testsheet2" -> "GUI"
# - Populate the (previously empty) API sheet with the same Sign-Up test
#   data/formatting that already lives on the GUI sheet
# - Update the GUI sheet's current selection
# - Update the API sheet's view (scrolled/selected near the new data)

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("testsheet1").Name = "API"
$wb.Worksheets.Item("testsheet2").Name = "GUI"

$apiWs = $wb.Worksheets.Item("API")
$guiWs = $wb.Worksheets.Item("GUI")

# Copy GUI's used range (values + formatting) onto API, starting at A1
$guiWs.UsedRange.Copy($apiWs.Range("A1"))

# The source row 2 only has data through column F (no G2 value) - the plain
# rectangular copy leaves a blank G2 behind, drop it so API matches GUI's shape
$apiWs.Range("G2").ClearContents()

# Give the columns on the new API sheet the same sort of best-fit custom
# widths that GUI already has
$apiWs.Columns.Item(1).ColumnWidth = 38.625
$apiWs.Columns.Item(2).ColumnWidth = 8.625
$apiWs.Columns.Item(3).ColumnWidth = 9.375
$apiWs.Columns.Item(4).ColumnWidth = 13.25
$apiWs.Columns.Item(5).ColumnWidth = 18.5
$apiWs.Columns.Item(6).ColumnWidth = 11.75
$apiWs.Columns.Item(7).ColumnWidth = 42.375

# Put the API sheet's view near the newly-added data (scrolled right one
# column, with G3 selected)
$apiWs.Activate()
$apiWs.Application.ActiveWindow.ScrollColumn = 2
$apiWs.Range("G3").Select()

# Update the GUI sheet's selection to the new-test-cases block, and leave it
# as the active/visible tab (matches the unchanged tabSelected="1")
$guiWs.Activate()
$guiWs.Range("A3:G4").Select()
